$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Czech sheet: insert a new row 14 ("FAT-S") above the existing "Wg"/
# "Accessories" rows, shifting them down to rows 15/16, and update the
# dimension + selection accordingly.
# ---------------------------------------------------------------------------
$wsCzech = $wb.Worksheets.Item("Czech")
$wsCzech.Rows("14:14").Insert()
$wsCzech.Range("A15").Copy()
$wsCzech.Range("A14").PasteSpecial(-4122)
$wsCzech.Range("A14").Value = "FAT-S"
$wsCzech.Range("A14").Select()

# ---------------------------------------------------------------------------
# Slovakia sheet: same row-insert pattern as Czech.
# ---------------------------------------------------------------------------
$wsSlovakia = $wb.Worksheets.Item("Slovakia")
$wsSlovakia.Rows("14:14").Insert()
$wsSlovakia.Range("A15").Copy()
$wsSlovakia.Range("A14").PasteSpecial(-4122)
$wsSlovakia.Range("A14").Value = "FAT-S"

# ---------------------------------------------------------------------------
# Update per-sheet selections/scroll state to match the saved view.
# ---------------------------------------------------------------------------
$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Range("A14").Select()

$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsBelgium.Rows("14:14").Select()

$wsSwiss = $wb.Worksheets.Item("Swiss")
$wsSwiss.Range("A12").Select()

$wsTurkey = $wb.Worksheets.Item("Turkey")
$wsTurkey.Range("H14").Select()

# Slovakia is the active/last-selected tab when the workbook was saved, so
# activate it last with its final selection.
$wsSlovakia.Range("A14").Select()
